$wb = $excel.ActiveWorkbook

# The order file encodes stimulus distance (D64/D80/D51) and size (S20/S30/S25)
# codes inside plain-text cell values (Condition, Filename_Left, Filename_Right,
# Distance) throughout the sheet. This regenerated order shifts the distances
# and the "large" size code to new values. Apply it as a global text
# substitution (whole workbook, every worksheet) using Excel's Find & Replace,
# with LookAt:=xlPart (2) so it rewrites the token wherever it appears inside
# a larger string (e.g. "Face14_D64_S20" -> "Face14_D69_S20",
# "Face09_D64_S30_l.png" -> "Face09_D69_S31_l.png").
#
# The four tokens are mutually disjoint (none of the replacement strings
# contains any of the other source tokens), so the substitutions are safe to
# run sequentially without clobbering each other.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("D64", "D69", 2)
    $ws.Cells.Replace("D80", "D86", 2)
    $ws.Cells.Replace("D51", "D55", 2)
    $ws.Cells.Replace("S30", "S31", 2)
}
